# fix test problem on importing xlsx
#
# Row 2 had the registrant's given name / family name swapped between
# columns A and B. Column A is "registrant_family_name" and column B is
# "registrant_given_name" (see row 1 headers), but A2 held "Wambui"
# (the given name) and B2 held "Njeri" (the family name). Swap them so
# A2 = "Njeri" (family name) and B2 = "Wambui" (given name).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$a2 = $ws.Range("A2").Value2
$b2 = $ws.Range("B2").Value2

$ws.Range("A2").Value2 = $b2
$ws.Range("B2").Value2 = $a2
